$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.591.05"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.922.96"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'246.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.70%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4722"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("D8").Value = "'0.2893"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.42%  "
$ws.Range("D9").Value = "'0.06790"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.38%  "
$ws.Range("D10").Value = "'105.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.00%  "
$ws.Range("D11").Value = "'18.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.920.09"
$ws.Range("E12").Value = "  +3.94%  "
$ws.Range("D13").Value = "'0.07699"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").Value = "'5.281"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.33%  "
$ws.Range("D15").Value = "'0.6727"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.97%  "
$ws.Range("D16").Value = "'291.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "30.594.39"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "'0.000007610"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").Value = "'0.9996"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'12.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").Value = "2.163.32"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D22").Value = "'5.464"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.80%  "
$ws.Range("D23").Value = "'0.9996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'6.333"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.23%  "
$ws.Range("D25").Value = "'9.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.02%  "
$ws.Range("D26").Value = "'168.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("D27").Value = "'20.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.25%  "
$ws.Range("E28").Value = "  +9.56%  "
$ws.Range("D29").Value = "'0.1076"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "'1.389"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.78%  "
$ws.Range("D31").Value = "'4.185"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").Value = "'4.130"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.31%  "
$ws.Range("D33").Value = "'0.05057"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").Value = "'0.7427"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").Value = "'1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.25%  "
$ws.Range("D36").Value = "'0.02084"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.46%  "
$ws.Range("D37").Value = "'2.748"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "'2.691"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").Value = "'2.052"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.56%  "
$ws.Range("D40").Value = "'111.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.41%  "
$ws.Range("D41").Value = "'0.8802"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.95%  "
$ws.Range("D42").Value = "'0.4401"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.39%  "
$ws.Range("D43").Value = "'5.893"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D44").Value = "'1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'67.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "'7.227"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").Value = "'9.344"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("D48").Value = "'47.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.20%  "
$ws.Range("D49").Value = "'0.1239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("D50").Value = "'34.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("D51").Value = "'0.4069"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.91%  "
